$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column CL: copy the width from column CK so the <col> metadata matches
$ws.Range("CL1").ColumnWidth = $ws.Range("CK1").ColumnWidth

# Header cell CL1: date text "2024/12/07" matching the style of CK1 (メイリオ font, no fill)
$ws.Range("CL1").NumberFormat = "@"
$ws.Range("CL1").Value = "2024/12/07"
$ws.Range("CK1").Copy()
$ws.Range("CL1").PasteSpecial(-4122)

# Data cells CL2:CL53 - values with conditional-style fills matching the source diff
$ws.Range("A2").Copy()
$ws.Range("CL2").PasteSpecial(-4122)
$ws.Range("CL2").Value = 153
$ws.Range("N2").Copy()
$ws.Range("CL3").PasteSpecial(-4122)
$ws.Range("CL3").Value = 128
$ws.Range("A2").Copy()
$ws.Range("CL4").PasteSpecial(-4122)
$ws.Range("CL4").Value = 160.4
$ws.Range("N2").Copy()
$ws.Range("CL5").PasteSpecial(-4122)
$ws.Range("CL5").Value = 133.8
$ws.Range("A2").Copy()
$ws.Range("CL6").PasteSpecial(-4122)
$ws.Range("CL6").Value = 170.7
$ws.Range("D2").Copy()
$ws.Range("CL7").PasteSpecial(-4122)
$ws.Range("CL7").Value = 112
$ws.Range("A2").Copy()
$ws.Range("CL8").PasteSpecial(-4122)
$ws.Range("CL8").Value = 155.6
$ws.Range("A2").Copy()
$ws.Range("CL9").PasteSpecial(-4122)
$ws.Range("CL9").Value = 178.3
$ws.Range("A2").Copy()
$ws.Range("CL10").PasteSpecial(-4122)
$ws.Range("CL10").Value = 158.4
$ws.Range("N2").Copy()
$ws.Range("CL11").PasteSpecial(-4122)
$ws.Range("CL11").Value = 129.8
$ws.Range("A2").Copy()
$ws.Range("CL12").PasteSpecial(-4122)
$ws.Range("CL12").Value = 174.4
$ws.Range("A2").Copy()
$ws.Range("CL13").PasteSpecial(-4122)
$ws.Range("CL13").Value = 165.9
$ws.Range("A2").Copy()
$ws.Range("CL14").PasteSpecial(-4122)
$ws.Range("CL14").Value = 142.9
$ws.Range("A2").Copy()
$ws.Range("CL15").PasteSpecial(-4122)
$ws.Range("CL15").Value = 163.4
$ws.Range("A2").Copy()
$ws.Range("CL16").PasteSpecial(-4122)
$ws.Range("CL16").Value = 190.2
$ws.Range("A2").Copy()
$ws.Range("CL17").PasteSpecial(-4122)
$ws.Range("CL17").Value = 141.5
$ws.Range("A2").Copy()
$ws.Range("CL18").PasteSpecial(-4122)
$ws.Range("CL18").Value = 145.6
$ws.Range("A2").Copy()
$ws.Range("CL19").PasteSpecial(-4122)
$ws.Range("CL19").Value = 150.9
$ws.Range("A2").Copy()
$ws.Range("CL20").PasteSpecial(-4122)
$ws.Range("CL20").Value = 159.9
$ws.Range("A2").Copy()
$ws.Range("CL21").PasteSpecial(-4122)
$ws.Range("CL21").Value = 170.7
$ws.Range("N2").Copy()
$ws.Range("CL22").PasteSpecial(-4122)
$ws.Range("CL22").Value = 132.5
$ws.Range("N2").Copy()
$ws.Range("CL23").PasteSpecial(-4122)
$ws.Range("CL23").Value = 127.6
$ws.Range("A2").Copy()
$ws.Range("CL24").PasteSpecial(-4122)
$ws.Range("CL24").Value = 160.7
$ws.Range("N2").Copy()
$ws.Range("CL25").PasteSpecial(-4122)
$ws.Range("CL25").Value = 126.6
$ws.Range("A2").Copy()
$ws.Range("CL26").PasteSpecial(-4122)
$ws.Range("CL26").Value = 145
$ws.Range("A2").Copy()
$ws.Range("CL27").PasteSpecial(-4122)
$ws.Range("CL27").Value = 149.4
$ws.Range("D2").Copy()
$ws.Range("CL28").PasteSpecial(-4122)
$ws.Range("CL28").Value = 119.6
$ws.Range("N2").Copy()
$ws.Range("CL29").PasteSpecial(-4122)
$ws.Range("CL29").Value = 139.6
$ws.Range("N2").Copy()
$ws.Range("CL30").PasteSpecial(-4122)
$ws.Range("CL30").Value = 131.6
$ws.Range("A2").Copy()
$ws.Range("CL31").PasteSpecial(-4122)
$ws.Range("CL31").Value = 145.3
$ws.Range("D2").Copy()
$ws.Range("CL32").PasteSpecial(-4122)
$ws.Range("CL32").Value = 123.4
$ws.Range("A2").Copy()
$ws.Range("CL33").PasteSpecial(-4122)
$ws.Range("CL33").Value = 140
$ws.Range("A2").Copy()
$ws.Range("CL34").PasteSpecial(-4122)
$ws.Range("CL34").Value = 163.9
$ws.Range("A2").Copy()
$ws.Range("CL35").PasteSpecial(-4122)
$ws.Range("CL35").Value = 145.4
$ws.Range("N2").Copy()
$ws.Range("CL36").PasteSpecial(-4122)
$ws.Range("CL36").Value = 135.1
$ws.Range("A2").Copy()
$ws.Range("CL37").PasteSpecial(-4122)
$ws.Range("CL37").Value = 145.6
$ws.Range("A2").Copy()
$ws.Range("CL38").PasteSpecial(-4122)
$ws.Range("CL38").Value = 151.6
$ws.Range("A2").Copy()
$ws.Range("CL39").PasteSpecial(-4122)
$ws.Range("CL39").Value = 160
$ws.Range("N2").Copy()
$ws.Range("CL40").PasteSpecial(-4122)
$ws.Range("CL40").Value = 129.4
$ws.Range("A2").Copy()
$ws.Range("CL41").PasteSpecial(-4122)
$ws.Range("CL41").Value = 157.4
$ws.Range("A2").Copy()
$ws.Range("CL42").PasteSpecial(-4122)
$ws.Range("CL42").Value = 178.6
$ws.Range("A2").Copy()
$ws.Range("CL43").PasteSpecial(-4122)
$ws.Range("CL43").Value = 186
$ws.Range("A2").Copy()
$ws.Range("CL44").PasteSpecial(-4122)
$ws.Range("CL44").Value = 144.5
$ws.Range("A2").Copy()
$ws.Range("CL45").PasteSpecial(-4122)
$ws.Range("CL45").Value = 145.7
$ws.Range("N2").Copy()
$ws.Range("CL46").PasteSpecial(-4122)
$ws.Range("CL46").Value = 132.5
$ws.Range("D2").Copy()
$ws.Range("CL47").PasteSpecial(-4122)
$ws.Range("CL47").Value = 123.5
$ws.Range("A2").Copy()
$ws.Range("CL48").PasteSpecial(-4122)
$ws.Range("CL48").Value = 174.6
$ws.Range("A2").Copy()
$ws.Range("CL49").PasteSpecial(-4122)
$ws.Range("CL49").Value = 187.5
$ws.Range("A2").Copy()
$ws.Range("CL50").PasteSpecial(-4122)
$ws.Range("CL50").Value = 150.8
$ws.Range("A2").Copy()
$ws.Range("CL51").PasteSpecial(-4122)
$ws.Range("CL51").Value = 225.3
$ws.Range("A2").Copy()
$ws.Range("CL52").PasteSpecial(-4122)
$ws.Range("CL52").Value = 151.1
$ws.Range("A2").Copy()
$ws.Range("CL53").PasteSpecial(-4122)
$ws.Range("CL53").Value = 144.6

$excel.CutCopyMode = 0